$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Headers
$ws.Range("D1").Value = "ep_icd10"
$ws.Range("E1").Value = "Specialty"

# Diagnosis (D) values - only present at the start of each diagnosis grouping
$ws.Range("D3").Value  = "A - head hurts"
$ws.Range("D5").Value  = "B - knee pain"
$ws.Range("D17").Value = "A - appendicitis"
$ws.Range("D27").Value = "A - fractures in several areas"
$ws.Range("D30").Value = "B- Sepsis"

# Specialty (E) values
$ws.Range("E3").Value  = "Neurology"
$ws.Range("E5").Value  = "Orthopaedics"
$ws.Range("E9").Value  = "Orthopaedics"
$ws.Range("E10").Value = "Cardiology"
$ws.Range("E11").Value = "Cardiology"
$ws.Range("E13").Value = "Gastroenterology"
$ws.Range("E14").Value = "Gastroenterology"
$ws.Range("E16").Value = "Cardiology"
$ws.Range("E17").Value = "Gastroenterology"
$ws.Range("E21").Value = "Gastroenterology"
$ws.Range("E22").Value = "Opthamology"
$ws.Range("E23").Value = "Opthamology"
$ws.Range("E24").Value = "Opthamology"
$ws.Range("E25").Value = "Opthamology"
$ws.Range("E27").Value = "Orthopaedics"
$ws.Range("E30").Value = "General Surgery"
$ws.Range("E34").Value = "General Surgery"
$ws.Range("E35").Value = "Opthamology"
$ws.Range("E36").Value = "Neurology"

$ws.Range("E31").Select()
